$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header value in column E (will also extend the shared strings table)
$ws.Range("E1").Value = "product-image"

# Set the column width for the new column E to match the diff (stored width 24).
# Excel quantizes ColumnWidth (character units) to whole pixels before storing the
# "width" attribute, so 23.17 is chosen as it reliably rounds to a stored width of 24.
$ws.Columns.Item(5).ColumnWidth = 23.17

# Update the selected cell to match the diff (F2)
$ws.Range("F2").Select()
